# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This updates column G ("K") on Sheet1 with freshly calculated values for
# each observation row (rows 2-67), replacing the previously stored
# "Strike#" derived figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values keyed by worksheet row number (row 2 = first data record).
$kValues = [ordered]@{
    2  = 0
    3  = 0
    4  = 1
    5  = 1
    6  = 1
    7  = 1
    8  = 0
    9  = 1
    10 = 1
    11 = 1
    12 = 1
    13 = 1
    14 = 0
    15 = 0
    16 = 0
    17 = 0
    18 = 0
    19 = 1
    20 = 3
    21 = 1
    22 = 0
    23 = 1
    24 = 3
    25 = 2
    26 = 3
    27 = 1
    28 = 1
    29 = 0
    30 = 1
    31 = 2
    32 = 3
    33 = 1
    34 = 0
    35 = 1
    36 = 1
    37 = 1
    38 = 2
    39 = 0
    40 = 0
    41 = 0
    42 = 0
    43 = 4
    44 = 2
    45 = 4
    46 = 1
    47 = 0
    48 = 1
    49 = 0
    50 = 2
    51 = 1
    52 = 0
    53 = 0
    54 = 2
    55 = 1
    56 = 1
    57 = 0
    58 = 1
    59 = 0
    60 = 0
    61 = 2
    62 = 3
    63 = 0
    64 = 0
    65 = 2
    66 = 2
    67 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
